# OptiStyle workbook update:
#  - Navigation sheet: add IconClass/Target columns, populate for Home &
#    Products rows, and append a new "Book Now" menu item.
#  - Products sheet: insert a new "Laser Beam Superman Shades" product row.
#  - Restore the selection/active-sheet state captured in the saved file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Navigation sheet
# ---------------------------------------------------------------
$wsNav = $wb.Worksheets.Item("Navigation")

$wsNav.Range("C1").Value = "IconClass"
$wsNav.Range("D1").Value = "Target"

$wsNav.Range("C2").Value = "fas fa-home"
$wsNav.Range("D2").Value = "_self"

$wsNav.Range("C3").Value = "fas fa-glasses"
$wsNav.Range("D3").Value = "_blank"

$wsNav.Range("A6").Value = "Book Now"
$wsNav.Range("B6").Value = "#appointment"

# ---------------------------------------------------------------
# Products sheet - insert a new row 5 for the new product
# ---------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("Products")

$wsProd.Rows.Item(5).Insert()
$wsProd.Range("A5").Value = "Laser Beam Superman Shades"
$wsProd.Range("B5").Value = "Fire Laserbeams"
$wsProd.Range("C5").Value = 9999.99
$wsProd.Range("C5").NumberFormat = "General"
$wsProd.Range("D5").Value = "fas fa-glasses text-red-600"

# ---------------------------------------------------------------
# Restore view state: Products selection, then Navigation active
# ---------------------------------------------------------------
$wsProd.Range("D13").Select()

$wsNav.Activate()
$wsNav.Range("H11").Select()
